$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.9
$ws.Range("H2").Value = 3.9
$ws.Range("I2").Value = 1.8
$ws.Range("J2").Value = 4.5
$ws.Range("L2").Value = 2.4
$ws.Range("U2").Value = 1.73
$ws.Range("V2").Value = 2
$ws.Range("W2").Value = 13
$ws.Range("AH2").Value = 8
$ws.Range("AI2").Value = 9
$ws.Range("AP2").Value = 29
$ws.Range("AX2").Value = 9.5
$ws.Range("AZ2").Value = 29

# Row 3
$ws.Range("G3").Value = 1.2
$ws.Range("H3").Value = 6.5
$ws.Range("I3").Value = 12
$ws.Range("J3").Value = 1.57
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 9
$ws.Range("M3").Value = 1.02
$ws.Range("N3").Value = 21
$ws.Range("O3").Value = 1.11
$ws.Range("P3").Value = 6.5
$ws.Range("Q3").Value = 1.36
$ws.Range("R3").Value = 3.1
$ws.Range("S3").Value = 1.2
$ws.Range("T3").Value = 4.33
$ws.Range("U3").Value = 1.83
$ws.Range("V3").Value = 1.83
$ws.Range("W3").Value = 10
$ws.Range("X3").Value = 7.5
$ws.Range("Y3").Value = 10
$ws.Range("Z3").Value = 8
$ws.Range("AA3").Value = 10
$ws.Range("AB3").Value = 23
$ws.Range("AC3").Value = 21
$ws.Range("AD3").Value = 13
$ws.Range("AG3").Value = 251
$ws.Range("AH3").Value = 34
$ws.Range("AI3").Value = 51
$ws.Range("AJ3").Value = 34
$ws.Range("AK3").Value = 151
$ws.Range("AL3").Value = 67
$ws.Range("AO3").Value = 5.5
$ws.Range("AP3").Value = 15
$ws.Range("AQ3").Value = 12
$ws.Range("AR3").Value = 29
$ws.Range("AT3").Value = 4.33
$ws.Range("AU3").Value = 9.5
$ws.Range("AW3").Value = 12
$ws.Range("AZ3").Value = 201
$ws.Range("BB3").Value = 251

# Row 4
$ws.Range("G4").Value = 1.98
$ws.Range("H4").Value = 3.45
$ws.Range("I4").Value = 3.4
$ws.Range("J4").Value = 2.55
$ws.Range("K4").Value = 2.15
$ws.Range("L4").Value = 3.75
$ws.Range("O4").Value = 1.24
$ws.Range("P4").Value = 3.3
$ws.Range("Q4").Value = 1.72
$ws.Range("R4").Value = 1.9
$ws.Range("U4").Value = 1.62
$ws.Range("V4").Value = 2.02
$ws.Range("W4").Value = 8
$ws.Range("X4").Value = 10
$ws.Range("Y4").Value = 8.5
$ws.Range("Z4").Value = 18
$ws.Range("AA4").Value = 15
$ws.Range("AB4").Value = 24
$ws.Range("AC4").Value = 11.5
$ws.Range("AD4").Value = 6.8
$ws.Range("AE4").Value = 13
$ws.Range("AF4").Value = 55
$ws.Range("AH4").Value = 11.75
$ws.Range("AI4").Value = 19.5
$ws.Range("AJ4").Value = 11.5
$ws.Range("AK4").Value = 45
$ws.Range("AL4").Value = 28
$ws.Range("AM4").Value = 32
$ws.Range("AN4").Value = 3.95
$ws.Range("AO4").Value = 10
$ws.Range("AP4").Value = 17.5
$ws.Range("AQ4").Value = 37
$ws.Range("AR4").Value = 65
$ws.Range("AT4").Value = 2.85
$ws.Range("AU4").Value = 6.8
$ws.Range("AW4").Value = 5.3
$ws.Range("AX4").Value = 18
$ws.Range("AY4").Value = 23
$ws.Range("AZ4").Value = 90
$ws.Range("BA4").Value = 110

# Row 6
$ws.Range("BC6").Value = 126
$ws.Range("BD6").Value = 126

# Row 9
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 3.2

# Row 10
$ws.Range("N10").Value = 12
$ws.Range("U10").Value = 1.73
$ws.Range("V10").Value = 2
$ws.Range("W10").Value = 8
$ws.Range("X10").Value = 8.5
$ws.Range("AC10").Value = 12
$ws.Range("AE10").Value = 15
$ws.Range("AG10").Value = 201
$ws.Range("BA10").Value = 101
$ws.Range("BB10").Value = 201
